# Y4_B2526_General_&_special_internal_1_reference_data.xlsx refresh
#
# The reference-data export was regenerated: the "Source File" column
# (column E) now records the name of this workbook itself (with a new
# generation timestamp) instead of the old "Group B1 2025-2026.xlsx"
# source name, column E was widened to fit the longer text, and the
# worksheet's page margins were reset to Excel's standard defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Source File) from 25 -> 50 characters so the longer
# file name fits without truncation.
$ws.Columns.Item(5).ColumnWidth = 49.17

# Reset page margins to Excel's built-in defaults (inches -> points,
# 1 inch = 72 points): left/right 0.75", top/bottom 1", header/footer 0.5".
$ws.PageSetup.LeftMargin   = 54
$ws.PageSetup.RightMargin  = 54
$ws.PageSetup.TopMargin    = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Update the "Source File" value (column E) for every data row (2-307):
# the export now stamps its own (renamed) file name instead of the old
# "Group B1 2025-2026.xlsx" source name.
$newSourceFile = "Y4_B2526_General_&_special_internal_1_reference_data_D07092025T122547.xlsx"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 307) { $lastRow = 307 }
$ws.Range("E2:E" + $lastRow).Value = $newSourceFile
